$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Rename the first sheet ---
$ws1.Name = "Primera prueba"

# --- Add the second sheet right after the first one ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "segunda prueba"

# =====================================================================
# Build "segunda prueba" content (mirrors "Primera prueba" layout/style)
# =====================================================================

# Row 1 - headers, copy formatting from Primera prueba row 1
$ws1.Range("A1:B1").Copy($ws2.Range("A1"))
$ws1.Range("E1:H1").Copy($ws2.Range("E1"))
$ws2.Range("A1").Value = "Comandos"
$ws2.Range("B1").Value = "Resultado"
$ws2.Range("E1").Value = "Evolución de Estados"

# Row 2
$ws1.Range("E2:H2").Copy($ws2.Range("E2"))
$ws2.Range("A2").Value = "write hal hobj"
$ws2.Range("B2").Value = "BAD_INSTRUCCION"
$ws2.Range("E2").Value = "Sujetos"
$ws2.Range("G2").Value = "Objetos"

# Row 3
$ws2.Range("A3").Value = "read hal"
$ws2.Range("B3").Value = "BAD_INSTRUCCION"
$ws2.Range("E3").Value = "lyle"
$ws2.Range("F3").Value = "hal"
$ws2.Range("G3").Value = "lobj"
$ws2.Range("H3").Value = "hobj"

# Row 4
$ws1.Range("E5").Copy($ws2.Range("E4"))
$ws1.Range("G5").Copy($ws2.Range("G4"))
$ws1.Range("E5").Copy($ws2.Range("H4"))
$ws2.Range("A4").Value = "write lyle lobj 10"
$ws2.Range("B4").Value = "OK"
$ws2.Range("E4").Value = 0
$ws2.Range("F4").Value = 10
$ws2.Range("G4").Value = 10
$ws2.Range("H4").Value = 20

# Row 5
$ws1.Range("E5").Copy($ws2.Range("F5"))
$ws2.Range("A5").Value = "read hal lobj"
$ws2.Range("B5").Value = "OK"
$ws2.Range("F5").Value = 20

# Row 6
$ws2.Range("A6").Value = "write lyle hobj 20"
$ws2.Range("B6").Value = "OK"

# Row 7
$ws2.Range("A7").Value = "write hal lobj 200"
$ws2.Range("B7").Value = "BAD_INSTRUCCION"

# Row 8
$ws2.Range("A8").Value = "read hal hobj"
$ws2.Range("B8").Value = "OK"

# Row 9
$ws2.Range("A9").Value = "read lyle lobj"
$ws2.Range("B9").Value = "OK"

# Row 10
$ws2.Range("A10").Value = "read lyle hobj"
$ws2.Range("B10").Value = "PROBLEMA PERMISOS"

# Row 11
$ws2.Range("A11").Value = "foo lyle lobj"
$ws2.Range("B11").Value = "BAD_INSTRUCCION"

# Row 12
$ws2.Range("A12").Value = "Hi lyle, This is hal"
$ws2.Range("B12").Value = "BAD_INSTRUCCION"

# Row 13
$ws2.Range("A13").Value = "The missile launch code is 1234567"
$ws2.Range("B13").Value = "BAD_INSTRUCCION"

# Merged cells on "segunda prueba"
$ws2.Range("E1:H1").Merge()
$ws2.Range("E2:F2").Merge()
$ws2.Range("G2:H2").Merge()

# Column widths / hidden columns on "segunda prueba"
$ws2.Columns.Item(1).ColumnWidth = 19.0
$ws2.Columns.Item(2).ColumnWidth = 21.75
$ws2.Columns.Item(3).ColumnWidth = 0
$ws2.Columns.Item(3).Hidden = $true
$ws2.Columns.Item(4).ColumnWidth = 0
$ws2.Columns.Item(4).Hidden = $true

# Page setup for both sheets
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# =====================================================================
# Update "Primera prueba": E5 value and selection / view state
# =====================================================================
$ws1.Range("E5").Value = 0

$ws1.Activate()
$ws1.Range("E5").Select()

$ws2.Activate()
$ws2.Range("E12").Select()
